$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.403.47"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "1.563.46"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'207.92"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'21.81"
$ws.Range("E8").Value = "  -2.25%  "
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.786.72"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("D13").Value = "1.573.82"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "'3.81"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "'0.513"
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").Value = "'63.29"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "27.406.16"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'212.01"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'7.24"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'4.10"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "'9.51"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Value = "'2.00"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").Value = "'153.22"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'6.70"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "'14.96"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -2.04%  "
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").Value = "1.361.02"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").Value = "'0.529"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E40").Value = "  +0.76%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'64.00"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "1.699.04"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("D48").Value = "'85.46"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "'0.0953"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("E51").Value = "  -0.83%  "
